$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "2025-09-01"
$ws.Range("A17").Style = "Normal"
$ws.Range("B17").Value = 57.86000061035156
$ws.Range("C17").Value = 690.1500244140625
$ws.Range("D17").Value = 321.1000061035156
